# "Cambios de estilos en home" - add a new student/row to the schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new row's values first, in the same order the target workbook's
# shared-string table was built (schedule text "Sabado  8:30 12:31" before the
# student name "ALUMNO, Alumno"), so new <si> entries land in the same order.
$ws.Range("C30").Value = "Sabado  8:30 12:31"
$ws.Range("B30").Value = "ALUMNO, Alumno"
$ws.Range("A30").Value = 121212

# Copy the formatting used by the existing rows onto the new row so the new
# cells keep the same borders/fonts/number formats as the rest of the table.
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C29").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C29").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Re-apply the values (PasteSpecial of formats shouldn't touch them, but make
# sure nothing got clobbered).
$ws.Range("A30").Value = 121212
$ws.Range("B30").Value = "ALUMNO, Alumno"
$ws.Range("C30").Value = "Sabado  8:30 12:31"

# Match the author's final selection state.
$ws.Range("H25").Select() | Out-Null

$wb.Save()
